# Micromouse code without delays.
#
# Renames Sheet1 -> "IR Sensor to CM Conversion" and Sheet2 ->
# "Tuple Recognition Statistics", fills the second sheet with the tuple /
# click-count recognition statistics table (and its helper notes), and
# leaves that sheet as the active tab/selection, matching the authored
# change.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Name = "IR Sensor to CM Conversion"
$ws2.Name = "Tuple Recognition Statistics"

# ---------------------------------------------------------------------
# Numeric source data + formulas for rows 3-16 (Units/L clicks/R clicks,
# per-unit approximations, averages)
# ---------------------------------------------------------------------
$aVals = @(1,    1,    1,    2,    2,    2,    3,    3,    3,    4,    4,    4,    4,    11)
$bVals = @(21,   28,   22,   1002, 1009, 1104, 2073, 2008, 1994, 3038, 3072, 3023, 3064, 10110)
$cVals = @(990,  998,  1005, 1949, 1929, 1988, 2988, 2969, 2909, 4033, 3998, 4006, 3937, 10973)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $r = 3 + $i
    $ws2.Range("A$r").Value = $aVals[$i]
    $ws2.Range("B$r").Value = $bVals[$i]
    $ws2.Range("C$r").Value = $cVals[$i]
    $ws2.Range("D$r").Formula = "=C$r/A$r"
    $ws2.Range("H$r").Formula = "=B$r/A$r"
    $ws2.Range("I$r").Formula = "=(H$r+D$r)/2"
}

# rows 17-20: only "Units Traveled" (A) recorded, rest derived from blanks
for ($r = 17; $r -le 20; $r++) {
    $ws2.Range("A$r").Value = 4
    $ws2.Range("D$r").Formula = "=C$r/A$r"
    $ws2.Range("H$r").Formula = "=B$r/A$r"
    $ws2.Range("I$r").Formula = "=(H$r+D$r)/2"
}

$ws2.Range("E3").Formula = "=AVERAGE(D3:D16)"

# ---------------------------------------------------------------------
# Title + header labels and free-form notes (written in the same order
# they were authored so shared-string ids line up)
# ---------------------------------------------------------------------
$ws2.Range("A2").Value = "Units Traveled"
$ws2.Range("B2").Value = "L clicks"
$ws2.Range("C2").Value = "R clicks"
$ws2.Range("H2").Value = "Approx 1 Unit L"
$ws2.Range("D2").Value = "Approx 1 Unit R"
$ws2.Range("I2").Value = "Approx 1 Unit Average"
$ws2.Range("J2").Value = "Outliers L"
$ws2.Range("K2").Value = "Outliers R"
$ws2.Range("J15").Value = "<-- started reprogramming after every poll to figure out if debugger is causing the outliers"
$ws2.Range("M1").Value = "35 clock pulses per 50 ms"
$ws2.Range("M2").Value = "WE WILL USE TIMER1; TIMER0 DOESN'T WORK!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!"
$ws2.Range("E2").Value = "Average"

# ---------------------------------------------------------------------
# Title row formatting + merge
# ---------------------------------------------------------------------
$titleRng = $ws2.Range("A1:G1")
$titleRng.HorizontalAlignment = -4108
$titleRng.Font.Bold = $true
$titleRng.Font.Size = 14
$ws2.Range("A1").Value = "Tuple Recognition Statistics"
$titleRng.Merge() | Out-Null
$ws2.Rows.Item(1).RowHeight = 18.75

# ---------------------------------------------------------------------
# Column widths (best-fit values from the authored workbook)
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 13.85546875
$ws2.Columns.Item(4).ColumnWidth = 27.5703125
$ws2.Columns.Item(5).ColumnWidth = 14.42578125
$ws2.Columns.Item(6).ColumnWidth = 14.7109375
$ws2.Columns.Item(7).ColumnWidth = 21
$ws2.Columns.Item(8).ColumnWidth = 9.42578125
$ws2.Columns.Item(9).ColumnWidth = 9.7109375

$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Selection + active tab
# ---------------------------------------------------------------------
$ws2.Range("E4").Select() | Out-Null
$ws2.Activate() | Out-Null
